$d = $word.ActiveDocument

function Get-PkgXml([string]$innerParagraphXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerParagraphXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-RunsByText([string]$needle, [string]$innerParagraphXml) {
    # Locate the needle text anywhere in the document, then rebuild a plain
    # Range (so Find's internal state doesn't interfere) and swap its
    # contents for freshly split <w:r> runs via InsertXML.
    $probe = $d.Content.Duplicate
    $probe.Find.Execute($needle, $true) | Out-Null
    $target = $d.Range($probe.Start, $probe.End)
    $null = $target.InsertXML((Get-PkgXml $innerParagraphXml))
}

# --- Change (paragraph "{m:endlet}"): "{m:" / "endlet}" -> "{" / "m:" / "endlet" / "}" ---
Replace-RunsByText "{m:endlet}" (
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m:</w:t></w:r>' +
    '<w:r><w:t>endlet</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>'
)

# --- Change (paragraph "...{m:self},"): "{m" / ":self}" -> "{" / "m" / ":self" / "}" ---
Replace-RunsByText "{m:self}," (
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:self</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '<w:r><w:t>,</w:t></w:r>'
)

# --- Change (paragraph "{m:let self= self.name}"): insert a space run between " self" and "=" ---
Replace-RunsByText " self= self.name}" (
    '<w:r><w:t xml:space="preserve"> self</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>=</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> self.</w:t></w:r>' +
    '<w:r><w:t>name}</w:t></w:r>'
)


